$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the "Selection Sort" row (row 9) data, which was previously blank,
# matching the pattern of the other completed rows.
$ws.Range("C9").Value = 37500
$ws.Range("E9").Value = 15.9618
$ws.Range("G9").Value = 64.9102

# Update the active cell selection on the sheet.
$ws.Range("H20").Select()
